$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test-case rows 13-15: fill in the previously-empty cells -----------
# Row 13: "ARE_START_7" / autostart-per-command-line-parameter test case
$ws.Range("A13").Value = "ARE_START_7"
$ws.Range("B13").Value = "Define autostart model per command line parameter"
$ws.Range("C13").Value = "TestModelAutostart.acs`nARE start file: start.bat (start.sh - Linux, Mac OS X)`nCommand shell: cmd.exe (xterm or similar - Linux, Mac OS X)"
$ws.Range("D13").Value = "0. Copy model to bin/ARE/models`n1. Open command shell in bin/ARE directory`n2. Execute `nstart.bat TestModelAutostart.acs`n./start.sh TestModelAutostart.acs"
$ws.Range("E13").Value = "The ARE must start successfully and the model 'TestModelAutostart.acs' must be started automatically"

# Row 14: "ARE_START_8" / start ARE without webservice
$ws.Range("A14").Value = "ARE_START_8"
$ws.Range("B14").Value = "Start ARE without webservice"
$ws.Range("C14").Value = "ARE_START_7"
$ws.Range("D14").Value = "0. Copy model to bin/ARE/models`n1. Open http://localhost:8082/`n2. Open http://localhost:8081/rest/runtime/model`n"
$ws.Range("E14").Value = "1. The page must not be loadable (err_connection_refused)`n2. The page must not be loadable (err_connection_refused)"

# Row 15: "ARE_START_9" / start ARE with webservice
$ws.Range("A15").Value = "ARE_START_9"
$ws.Range("B15").Value = "Start ARE with webservice"
$ws.Range("C15").Value = "TestModelAutostart.acs`nARE start file: start.bat (start.sh - Linux, Mac OS X)`nCommand shell: cmd.exe (xterm or similar - Linux, Mac OS X)"
$ws.Range("D15").Value = "0. Copy model to bin/ARE/models`n1. Open command shell in bin/ARE directory`n2. Execute `nstart.bat --webservice TestModelAutostart.acs`n./start.sh --webservice TestModelAutostart.acs`n3. Open http://localhost:8082/`n4. Open http://localhost:8081/rest/runtime/model`n"
$ws.Range("E15").Value = "1. The page must be loadable and stating 'CONNECTED'  to the websocket`n2. The page must be loadable and return the currently deployed model in xml"

# --- Scroll position: the sheet view now shows row 9 at the top -------------
# (best-effort; restore the original selection afterwards so only the
# scroll/topLeftCell position changes)
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
[void]$excel.Goto($ws.Range("A9"), $true)
[void]$ws.Range("D13").Select()
